# Edit script: reorder the Odd_CS_* columns AB..AH (insert Odd_CS_0-1 before
# the old Odd_CS_4-4 column, pushing Odd_CS_4-4 to the end), and fill in the
# odds data for rows 2 and 5 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): columns AB1:AH1 ---
$headers = @("Odd_CS_0-1", "Odd_CS_0-2", "Odd_CS_1-2", "Odd_CS_0-3", "Odd_CS_1-3", "Odd_CS_2-3", "Odd_CS_4-4")
$startCol = 28  # column AB
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $startCol + $i).Value = $headers[$i]
}

# --- Row 2 (G2:AG2) now filled with odds values; AH2 remains empty ---
$row2 = @(2.1, 2.95, 3.65, 1.47, 2.32, 2.37, 1.45, 1.55, 2.15, 2.02, 1.62, 5.5, 8.75, 9.25, 19.5, 21, 40, 6.4, 5.9, 18, 120, 8.25, 18, 13, 55, 40, 55)
$row2StartCol = 7  # column G
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $row2StartCol + $i).Value = $row2[$i]
}

# --- Row 5 (AB5:AH5) reordered odds values ---
$row5 = @(8.75, 17.5, 12.5, 50, 37, 45, 800)
for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, $startCol + $i).Value = $row5[$i]
}
